$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.321.41"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").Value = "1.932.12"
$ws.Range("E3").Value = "  -0.10%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "'251.18"
$ws.Range("E5").Value = "  +1.24%  "
$ws.Range("D6").Value = "'0.7122"
$ws.Range("E6").Value = "  -1.38%  "
$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("D8").Value = "'0.3257"
$ws.Range("E8").Value = "  -1.22%  "
$ws.Range("D9").Value = "'27.31"
$ws.Range("E9").Value = "  +2.25%  "
$ws.Range("D10").Value = "'0.07189"
$ws.Range("E10").Value = "  +4.95%  "
$ws.Range("D11").Value = "'0.7980"
$ws.Range("E11").Value = "  -1.53%  "
$ws.Range("D12").Value = "'0.08093"
$ws.Range("E12").Value = "  +1.49%  "
$ws.Range("D13").Value = "1.929.52"
$ws.Range("E13").Value = "  -0.11%  "
$ws.Range("D14").Value = "'5.429"
$ws.Range("E14").Value = "  -0.42%  "
$ws.Range("D15").Value = "'94.80"
$ws.Range("E15").Value = "  -0.27%  "
$ws.Range("D16").Value = "'14.81"
$ws.Range("E16").Value = "  +1.31%  "
$ws.Range("D17").Value = "30.324.99"
$ws.Range("E17").Value = "  +0.12%  "
$ws.Range("D18").Value = "'252.44"
$ws.Range("E18").Value = "  -4.32%  "
$ws.Range("D19").Value = "'0.000008101"
$ws.Range("E19").Value = "  +1.52%  "
$ws.Range("D20").Value = "'5.792"
$ws.Range("E20").Value = "  -0.89%  "
$ws.Range("D21").Value = "2.182.13"
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("D23").Value = "'1.002"
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("D24").Value = "'6.914"
$ws.Range("E24").Value = "  -0.24%  "
$ws.Range("D25").Value = "'9.693"
$ws.Range("E25").Value = "  -0.56%  "
$ws.Range("D26").Value = "'164.76"
$ws.Range("E26").Value = "  +2.82%  "
$ws.Range("D27").Value = "'19.23"
$ws.Range("E27").Value = "  +0.72%  "
$ws.Range("D28").Value = "'2.314"
$ws.Range("E28").Value = "  -2.00%  "
$ws.Range("D29").Value = "'0.1280"
$ws.Range("E29").Value = "  -5.28%  "
$ws.Range("D30").Value = "'1.361"
$ws.Range("E30").Value = "  -0.53%  "
$ws.Range("D31").Value = "'1.543"
$ws.Range("E31").Value = "  -0.61%  "
$ws.Range("D32").Value = "'4.428"
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("D33").Value = "'4.198"
$ws.Range("E33").Value = "  -0.80%  "
$ws.Range("D34").Value = "'0.05202"
$ws.Range("E34").Value = "  +1.89%  "
$ws.Range("D35").Value = "'1.264"
$ws.Range("E35").Value = "  +4.27%  "
$ws.Range("D36").Value = "'0.7486"
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("D37").Value = "'2.768"
$ws.Range("E37").Value = "  +1.48%  "
$ws.Range("D38").Value = "'0.01962"
$ws.Range("E38").Value = "  +0.58%  "
$ws.Range("D39").Value = "'2.803"
$ws.Range("E39").Value = "  -0.76%  "
$ws.Range("D40").Value = "'78.98"
$ws.Range("E40").Value = "  -2.51%  "
$ws.Range("D41").Value = "'6.430"
$ws.Range("E41").Value = "  -2.44%  "
$ws.Range("D42").Value = "'0.4524"
$ws.Range("E42").Value = "  +0.49%  "
$ws.Range("D43").Value = "'2.025"
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("D44").Value = "'0.8420"
$ws.Range("E44").Value = "  +0.35%  "
$ws.Range("D45").Value = "'1.001"
$ws.Range("E45").Value = "  +0.13%  "
$ws.Range("D46").Value = "'101.86"
$ws.Range("E46").Value = "  -0.76%  "
$ws.Range("D47").Value = "'9.792"
$ws.Range("E47").Value = "  +0.30%  "
$ws.Range("D48").Value = "'7.431"
$ws.Range("E48").Value = "  +0.96%  "
$ws.Range("D49").Value = "'36.69"
$ws.Range("E49").Value = "  +0.74%  "
$ws.Range("D50").Value = "'0.06094"
$ws.Range("E50").Value = "  +2.70%  "
$ws.Range("D51").Value = "'0.4179"
$ws.Range("E51").Value = "  +0.92%  "
